# Refresh the cryptocurrency price/volume table with the latest values
# pulled on Tue Mar 21 01:17:13 UTC 2023 (GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.884.34"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.750.79"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "334.25"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.3781"
$ws.Range("E7").Value = "  -3.42%  "
$ws.Range("D8").Value = "'0.3354"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").Value = "'44.63"
$ws.Range("E9").Value = "  -7.40%  "
$ws.Range("D10").Value = "'1.114"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "'0.07192"
$ws.Range("E11").Value = "  -4.76%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'22.29"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'6.119"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("D15").Value = "'7.130"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "1.754.39"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'0.00001053"
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").Value = "'0.06573"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'78.74"
$ws.Range("E19").Value = "  -6.39%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "'16.70"
$ws.Range("E21").Value = "  -6.55%  "
$ws.Range("D22").Value = "'6.223"
$ws.Range("E22").Value = "  -5.36%  "
$ws.Range("D23").Value = "27.914.48"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "11.59"
$ws.Range("E24").Value = "  -6.56%  "
$ws.Range("D25").Value = "'2.394"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'152.67"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'19.83"
$ws.Range("E27").Value = "  -6.68%  "
$ws.Range("D28").Value = "'2.307"
$ws.Range("E28").Value = "  -8.11%  "
$ws.Range("D29").Value = "1.957.64"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.261"
$ws.Range("E30").Value = "  -17.38%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'131.40"
$ws.Range("E31").Value = "  -5.23%  "
$ws.Range("D32").Value = "'4.022"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'5.771"
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("D34").Value = "'0.08734"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "'12.11"
$ws.Range("E35").Value = "  -7.65%  "
$ws.Range("D36").Value = "'0.6651"
$ws.Range("E36").Value = "  -5.37%  "
$ws.Range("D37").Value = "'0.02304"
$ws.Range("E37").Value = "  -7.01%  "
$ws.Range("D38").Value = "'0.06202"
$ws.Range("E38").Value = "  -4.94%  "
$ws.Range("D39").Value = "'5.130"
$ws.Range("E39").Value = "  -6.94%  "
$ws.Range("D40").Value = "0.2102"
$ws.Range("E40").Value = "  -6.76%  "
$ws.Range("D41").Value = "'1.209"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("D42").Value = "1.462"
$ws.Range("E42").Value = "  -8.78%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.939"
$ws.Range("E43").Value = "  -7.31%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  -5.74%  "
$ws.Range("D46").Value = "'3.817"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'0.6010"
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("D48").Value = "126.03"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("D49").Value = "'2.005"
$ws.Range("E49").Value = "  -7.01%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.180"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'0.06994"
$ws.Range("E51").Value = "  -3.08%  "
